# Insert a new column before column A, shifting the existing "k:2"..."k:15"
# header row and the 10 data rows one column to the right (A1:N11 -> B1:O11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

# Fill the newly-created column A with fold labels for the 10 data rows.
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = "fold $i"
}

# Match the bold/centered/bordered header formatting (same style as row 1)
# by copying the format from a header cell onto the new label cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
